$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.486.39"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "'3.150.47"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'241.20"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "'619.17"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").Value = "'1.14"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").Value = "'0.377"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'3.149.19"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").Value = "'0.749"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "'5.62"
$ws.Range("E15").Value = "  +2.24%  "
$ws.Range("D16").Value = "'91.174.26"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("D17").Value = "'3.733.20"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").Value = "'3.155.37"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("D19").Value = "'3.75"
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("D20").Value = "'15.10"
$ws.Range("E20").Value = "  +5.41%  "
$ws.Range("E21").Value = "  +4.50%  "
$ws.Range("D22").Value = "'457.31"
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").Value = "'9.16"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").Value = "'5.94"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.64"
$ws.Range("E26").Value = "  +63.51%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "'89.21"
$ws.Range("E27").Value = "  -3.67%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'11.89"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.150"
$ws.Range("E30").Value = "  +32.92%  "
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.238"
$ws.Range("E32").Value = "  +11.11%  "
$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").Value = "'0.169"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'9.46"
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.177"
$ws.Range("E35").Value = "  +13.14%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'26.56"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'7.53"
$ws.Range("E37").Value = "  -3.85%  "
$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").Value = "'1.95"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'494.44"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "'1.33"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("B41").Value = "MantraDAO"
$ws.Range("C41").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D41").Value = "'3.88"
$ws.Range("E41").Value = "  -9.75%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.450"
$ws.Range("E42").Value = "  +7.70%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'3.43"
$ws.Range("E43").Value = "  -6.33%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "'22.14"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'0.716"
$ws.Range("E46").Value = "  +4.29%  "
$ws.Range("D47").Value = "'1.94"
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D48").Value = "'156.42"
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "'0.0329"
$ws.Range("E51").Value = "  +6.39%  "
